# Update countries & provincias Spain
# -------------------------------------------------------------------------
# paises.xlsx ("Pais" sheet) refresh: numbers updated for several countries
# and "Irak" / "Burkina Faso" re-ranked (each cut from its old spot and
# re-inserted a few rows earlier), which cascades a one-row shift onto the
# countries that used to sit between the old and new position. The net,
# cell-by-cell effect (label + 7 data columns B:H) is applied explicitly
# below, plus the "Datos actualizados..." timestamp in A1.
# -------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer / title timestamp.
$ws.Range("A1").Value = "Datos actualizados a 30 de Marzo de 2020 a las 19:20"

# Row => (new label if it changed, else $null) + new B..H values.
$data = @(
    @{Row=4;  Label=$null;                     B=146027; C=2536; D=4579; E=138798; F=3087; G=67; H=2650},
    @{Row=17; Label=$null;                     B=9541;   C=753;  D=636;  E=8797;   F=193;  G=22; H=108},
    @{Row=25; Label=$null;                     B=2942;   C=125;  D=25;   E=2894;   F=52;   G=7;  H=23},
    @{Row=32; Label=$null;                     B=1952;   C=137;  D=209;  E=1693;   F=33;   G=7;  H=50},
    @{Row=36; Label=$null;                     B=1690;   C=93;   D=76;   E=1593;   F=11;   G=7;  H=21},
    @{Row=60; Label="Irak";                     B=630;    C=83;   D=152;  E=432;    F=0;    G=4;  H=46},
    @{Row=61; Label="Emiratos Arabes Unidos";   B=611;    C=41;   D=61;   E=545;    F=2;    G=2;  H=5},
    @{Row=62; Label="Egipto";                   B=609;    C=0;    D=132;  E=437;    F=0;    G=0;  H=40},
    @{Row=63; Label="Nueva Zelanda";             B=589;    C=75;   D=63;   E=525;    F=2;    G=0;  H=1},
    @{Row=64; Label="Argelia";                   B=584;    C=73;   D=37;   E=512;    F=0;    G=4;  H=35},
    @{Row=65; Label=$null;                     B=534;    C=55;   D=14;   E=487;    F=1;    G=7;  H=33},
    @{Row=67; Label=$null;                     B=491;    C=31;   D=1;    E=483;    F=5;    G=0;  H=7},
    @{Row=87; Label="Burkina Faso";              B=246;    C=24;   D=31;   E=203;    F=0;    G=0;  H=12},
    @{Row=88; Label="Republica de Chipre";       B=230;    C=16;   D=15;   E=208;    F=3;    G=1;  H=7},
    @{Row=89; Label="San Marino";                 B=230;    C=6;    D=13;   E=192;    F=16;   G=3;  H=25},
    @{Row=90; Label="Reunion";                    B=224;    C=41;   D=1;    E=223;    F=0;    G=0;  H=0},
    @{Row=91; Label="Albania";                    B=223;    C=11;   D=44;   E=168;    F=7;    G=1;  H=11}
)

foreach ($item in $data) {
    $r = $item.Row
    if ($item.Label -ne $null) {
        $ws.Cells.Item($r, 1).Value = $item.Label
    }
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
    $ws.Cells.Item($r, 8).Value = $item.H
}
